# Horarios Línea 141 - refresh scrape (07:28:23) across the 3 sheets:
#   LP1912, LP1912-215, 6203-6173
# Mirrors the new scheduled-arrival rows appended by the scraper plus the
# handful of in-place corrections to already-printed rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2,1).Value = "Última actualización: 07:28:23"
$ws1.Cells.Item(3,1).Value = "Total filas: 59"

# Rows 22/23 swapped their "Linea" label
$ws1.Cells.Item(22,3).Value = "14_ABASTO"
$ws1.Cells.Item(23,3).Value = "215C_EL PATO"

# Rows 35-48: recomputed Hora_Scrap/Hora_Llegada/Linea/Minutos for the new scrape
$ws1.Cells.Item(35,1).Value = "07:28:23"
$ws1.Cells.Item(35,2).Value = "07:28"
$ws1.Cells.Item(35,4).Value = 0

$ws1.Cells.Item(36,1).Value = "07:28:23"
$ws1.Cells.Item(36,2).Value = "07:29"
$ws1.Cells.Item(36,4).Value = 1

$ws1.Cells.Item(37,2).Value = "07:31"
$ws1.Cells.Item(37,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(37,4).Value = 36

$ws1.Cells.Item(38,1).Value = "07:28:23"
$ws1.Cells.Item(38,2).Value = "07:31"
$ws1.Cells.Item(38,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(38,4).Value = 3

$ws1.Cells.Item(39,1).Value = "07:28:23"
$ws1.Cells.Item(39,2).Value = "07:32"
$ws1.Cells.Item(39,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(39,4).Value = 4

$ws1.Cells.Item(40,1).Value = "07:28:23"
$ws1.Cells.Item(40,2).Value = "07:34"
$ws1.Cells.Item(40,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(40,4).Value = 6

$ws1.Cells.Item(41,2).Value = "07:36"
$ws1.Cells.Item(41,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(41,4).Value = 41

$ws1.Cells.Item(42,1).Value = "07:28:23"
$ws1.Cells.Item(42,2).Value = "07:37"
$ws1.Cells.Item(42,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(42,4).Value = 9

$ws1.Cells.Item(43,1).Value = "07:28:23"
$ws1.Cells.Item(43,2).Value = "07:39"
$ws1.Cells.Item(43,3).Value = "10_OLMOS"
$ws1.Cells.Item(43,4).Value = 11

$ws1.Cells.Item(44,1).Value = "07:28:23"
$ws1.Cells.Item(44,2).Value = "07:47"
$ws1.Cells.Item(44,3).Value = "14_ABASTO"
$ws1.Cells.Item(44,4).Value = 19

$ws1.Cells.Item(45,1).Value = "07:28:23"
$ws1.Cells.Item(45,2).Value = "07:51"
$ws1.Cells.Item(45,3).Value = "215D_EL PATO"
$ws1.Cells.Item(45,4).Value = 23

$ws1.Cells.Item(46,1).Value = "07:28:23"
$ws1.Cells.Item(46,2).Value = "07:55"
$ws1.Cells.Item(46,3).Value = "10_OLMOS"
$ws1.Cells.Item(46,4).Value = 27

$ws1.Cells.Item(47,1).Value = "07:28:23"
$ws1.Cells.Item(47,2).Value = "08:03"
$ws1.Cells.Item(47,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(47,4).Value = 35

$ws1.Cells.Item(48,1).Value = "07:28:23"
$ws1.Cells.Item(48,2).Value = "08:08"
$ws1.Cells.Item(48,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(48,4).Value = 40

$ws1.Cells.Item(49,2).Value = "08:09"
$ws1.Cells.Item(49,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(49,4).Value = 74

# New rows 50-64 appended by the scrape
$ws1.Cells.Item(50,1).Value = "07:28:23"
$ws1.Cells.Item(50,2).Value = "08:12"
$ws1.Cells.Item(50,3).Value = "15_ABASTO"
$ws1.Cells.Item(50,4).Value = 44
$ws1.Cells.Item(50,5).Value = "LP1912"

$ws1.Cells.Item(51,1).Value = "07:28:23"
$ws1.Cells.Item(51,2).Value = "08:21"
$ws1.Cells.Item(51,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(51,4).Value = 53
$ws1.Cells.Item(51,5).Value = "LP1912"

$ws1.Cells.Item(52,1).Value = "07:28:23"
$ws1.Cells.Item(52,2).Value = "08:22"
$ws1.Cells.Item(52,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(52,4).Value = 54
$ws1.Cells.Item(52,5).Value = "LP1912"

$ws1.Cells.Item(53,1).Value = "07:28:23"
$ws1.Cells.Item(53,2).Value = "08:23"
$ws1.Cells.Item(53,3).Value = "215B_EL PATO"
$ws1.Cells.Item(53,4).Value = 55
$ws1.Cells.Item(53,5).Value = "LP1912"

$ws1.Cells.Item(54,1).Value = "07:28:23"
$ws1.Cells.Item(54,2).Value = "08:27"
$ws1.Cells.Item(54,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(54,4).Value = 59
$ws1.Cells.Item(54,5).Value = "LP1912"

$ws1.Cells.Item(55,1).Value = "06:55:48"
$ws1.Cells.Item(55,2).Value = "08:42"
$ws1.Cells.Item(55,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(55,4).Value = 107
$ws1.Cells.Item(55,5).Value = "LP1912"

$ws1.Cells.Item(56,1).Value = "07:28:23"
$ws1.Cells.Item(56,2).Value = "08:43"
$ws1.Cells.Item(56,3).Value = "14_ABASTO"
$ws1.Cells.Item(56,4).Value = 75
$ws1.Cells.Item(56,5).Value = "LP1912"

$ws1.Cells.Item(57,1).Value = "07:28:23"
$ws1.Cells.Item(57,2).Value = "08:50"
$ws1.Cells.Item(57,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(57,4).Value = 82
$ws1.Cells.Item(57,5).Value = "LP1912"

$ws1.Cells.Item(58,1).Value = "07:28:23"
$ws1.Cells.Item(58,2).Value = "08:54"
$ws1.Cells.Item(58,3).Value = "17_ROMERO"
$ws1.Cells.Item(58,4).Value = 86
$ws1.Cells.Item(58,5).Value = "LP1912"

$ws1.Cells.Item(59,1).Value = "07:28:23"
$ws1.Cells.Item(59,2).Value = "09:01"
$ws1.Cells.Item(59,3).Value = "215A_EL PATO"
$ws1.Cells.Item(59,4).Value = 93
$ws1.Cells.Item(59,5).Value = "LP1912"

$ws1.Cells.Item(60,1).Value = "07:28:23"
$ws1.Cells.Item(60,2).Value = "09:10"
$ws1.Cells.Item(60,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(60,4).Value = 102
$ws1.Cells.Item(60,5).Value = "LP1912"

$ws1.Cells.Item(61,1).Value = "07:28:23"
$ws1.Cells.Item(61,2).Value = "09:17"
$ws1.Cells.Item(61,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(61,4).Value = 109
$ws1.Cells.Item(61,5).Value = "LP1912"

$ws1.Cells.Item(62,1).Value = "07:28:23"
$ws1.Cells.Item(62,2).Value = "09:21"
$ws1.Cells.Item(62,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(62,4).Value = 113
$ws1.Cells.Item(62,5).Value = "LP1912"

$ws1.Cells.Item(63,1).Value = "07:28:23"
$ws1.Cells.Item(63,2).Value = "09:23"
$ws1.Cells.Item(63,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(63,4).Value = 115
$ws1.Cells.Item(63,5).Value = "LP1912"

$ws1.Cells.Item(64,1).Value = "07:28:23"
$ws1.Cells.Item(64,2).Value = "09:23"
$ws1.Cells.Item(64,3).Value = "17_ROMERO"
$ws1.Cells.Item(64,4).Value = 115
$ws1.Cells.Item(64,5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2,1).Value = "Última actualización: 07:28:23"
$ws2.Cells.Item(3,1).Value = "Total filas: 9"

$ws2.Cells.Item(12,1).Value = "07:28:23"
$ws2.Cells.Item(12,4).Value = 23

$ws2.Cells.Item(13,1).Value = "07:28:23"
$ws2.Cells.Item(13,4).Value = 55

# New row 14 appended by the scrape
$ws2.Cells.Item(14,1).Value = "07:28:23"
$ws2.Cells.Item(14,2).Value = "09:01"
$ws2.Cells.Item(14,3).Value = "215A_EL PATO"
$ws2.Cells.Item(14,4).Value = 93
$ws2.Cells.Item(14,5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2,1).Value = "Última actualización: 07:28:23"
$ws3.Cells.Item(3,1).Value = "Total filas: 14"

$ws3.Cells.Item(14,1).Value = "07:28:23"
$ws3.Cells.Item(14,4).Value = 7

# A brand-new row is inserted at 16, pushing the former row 16 down to 17
$ws3.Rows.Item(16).Insert()
$ws3.Cells.Item(16,1).Value = "07:28:23"
$ws3.Cells.Item(16,2).Value = "08:18"
$ws3.Cells.Item(16,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(16,4).Value = 50
$ws3.Cells.Item(16,5).Value = "L6203"

# New rows 18-19 appended by the scrape
$ws3.Cells.Item(18,1).Value = "07:28:23"
$ws3.Cells.Item(18,2).Value = "08:35"
$ws3.Cells.Item(18,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(18,4).Value = 67
$ws3.Cells.Item(18,5).Value = "L6173"

$ws3.Cells.Item(19,1).Value = "07:28:23"
$ws3.Cells.Item(19,2).Value = "09:09"
$ws3.Cells.Item(19,3).Value = "215D_LA PLATA"
$ws3.Cells.Item(19,4).Value = 101
$ws3.Cells.Item(19,5).Value = "L6203"
